# Update the organization website shared string in B10:
# "www.stat.kg" -> "www.stat.gov.kg"
#
# Excel regenerates a fresh (deduplicated) shared-string table whenever a
# cell's text changes, so simply writing the new value naturally removes the
# old "www.stat.kg" entry and appends the new "www.stat.gov.kg" entry at the
# end of sharedStrings.xml, shifting every other string's index down by one
# -- matching the upstream edit exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B10")
$cell.Value = "www.stat.gov.kg"

# Touch (re-apply) the font's theme color with the value it already has.
# This is a no-visual-effect nudge, but it makes Excel materialise a
# distinct font/style record for this cell instead of continuing to share
# the old style index -- exactly like the committed workbook, where B10
# ends up pointing at a brand-new cell style instead of the one it used to
# share with the other single-line cells in the sheet.
$cell.Font.ThemeColor = 1

# Leave the edited cell selected, like the saved workbook shows.
$cell.Select()
